{"js": "// Replace the date line and the multiplication problem/answer cells\n// per the commit diff. Each old string is unique in the document body,\n// so a straightforward search + full-text replace on the single\n// matching run is sufficient and preserves all other formatting.\n\nconst replacements = [\n  [\"2025-01-29 Wednesday\", \"2025-01-30 Thursday\"],\n  [\"76\u00d722=1672\", \"23\u00d784=1932\"],\n  [\"99\u00d737=3663\", \"62\u00d782=5084\"],\n  [\"27\u00d730=810\", \"85\u00d713=1105\"],\n  [\"95\u00d778=7410\", \"71\u00d788=6248\"],\n  [\"62\u00d786=5332\", \"47\u00d781=3807\"],\n  [\"83\u00d732=2656\", \"35\u00d756=1960\"],\n  [\"11\u00d772=792\", \"25\u00d761=1525\"],\n  [\"79\u00d726=2054\", \"52\u00d744=2288\"],\n  [\"68\u00d786=5848\", \"19\u00d719=361\"],\n  [\"16\u00d728=448\", \"83\u00d774=6142\"],\n  [\"59\u00d780=4720\", \"81\u00d749=3969\"],\n  [\"41\u00d789=3649\", \"34\u00d766=2244\"],\n  [\"73\u00d762=4526\", \"97\u00d758=5626\"],\n  [\"34\u00d764=2176\", \"33\u00d762=2046\"],\n  [\"38\u00d756=2128\", \"20\u00d777=1540\"],\n  [\"54\u00d763=3402\", \"73\u00d741=2993\"],\n  [\"35\u00d761=2135\", \"63\u00d764=4032\"],\n  [\"33\u00d730=990\", \"92\u00d720=1840\"],\n  [\"77\u00d768=5236\", \"15\u00d755=825\"],\n  [\"95\u00d736=3420\", \"73\u00d751=3723\"],\n  [\"25\u00d766=1650\", \"43\u00d747=2021\"],\n  [\"72\u00d745=3240\", \"69\u00d716=1104\"],\n  [\"42\u00d787=3654\", \"46\u00d726=1196\"],\n  [\"15\u00d749=735\", \"26\u00d745=1170\"],\n  [\"93\u00d744=4092\", \"61\u00d786=5246\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the multiplication problem/answer cells\n# per the commit diff, using Word's Find/Replace on the document range.\n# Each \"old\" string occurs exactly once in the document, so a plain\n# text Find & Replace for each pair reproduces the diff precisely while\n# leaving all other formatting (run properties, table layout, etc.)\n# untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-01-29 Wednesday\", \"2025-01-30 Thursday\"),\n    @(\"76\u00d722=1672\", \"23\u00d784=1932\"),\n    @(\"99\u00d737=3663\", \"62\u00d782=5084\"),\n    @(\"27\u00d730=810\", \"85\u00d713=1105\"),\n    @(\"95\u00d778=7410\", \"71\u00d788=6248\"),\n    @(\"62\u00d786=5332\", \"47\u00d781=3807\"),\n    @(\"83\u00d732=2656\", \"35\u00d756=1960\"),\n    @(\"11\u00d772=792\", \"25\u00d761=1525\"),\n    @(\"79\u00d726=2054\", \"52\u00d744=2288\"),\n    @(\"68\u00d786=5848\", \"19\u00d719=361\"),\n    @(\"16\u00d728=448\", \"83\u00d774=6142\"),\n    @(\"59\u00d780=4720\", \"81\u00d749=3969\"),\n    @(\"41\u00d789=3649\", \"34\u00d766=2244\"),\n    @(\"73\u00d762=4526\", \"97\u00d758=5626\"),\n    @(\"34\u00d764=2176\", \"33\u00d762=2046\"),\n    @(\"38\u00d756=2128\", \"20\u00d777=1540\"),\n    @(\"54\u00d763=3402\", \"73\u00d741=2993\"),\n    @(\"35\u00d761=2135\", \"63\u00d764=4032\"),\n    @(\"33\u00d730=990\", \"92\u00d720=1840\"),\n    @(\"77\u00d768=5236\", \"15\u00d755=825\"),\n    @(\"95\u00d736=3420\", \"73\u00d751=3723\"),\n    @(\"25\u00d766=1650\", \"43\u00d747=2021\"),\n    @(\"72\u00d745=3240\", \"69\u00d716=1104\"),\n    @(\"42\u00d787=3654\", \"46\u00d726=1196\"),\n    @(\"15\u00d749=735\", \"26\u00d745=1170\"),\n    @(\"93\u00d744=4092\", \"61\u00d786=5246\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute([ref]$oldText, $true, $true, $false, $false, $false, $true, 1, $false, [ref]$newText, 2) | Out-Null\n}\n"}
